$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.618.49"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "3.028.40"
$ws.Range("E3").Value = "  +4.48%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "201.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "630.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.31%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.210"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.37%  "

$ws.Range("D10").Value = "3.028.19"
$ws.Range("E10").Value = "  +4.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.437"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.51%  "

$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("E13").Value = "  +7.45%  "

$ws.Range("D14").Value = "3.586.33"
$ws.Range("E14").Value = "  +4.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.19%  "

$ws.Range("D16").Value = "76.528.40"
$ws.Range("E16").Value = "  +1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000194"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "3.028.38"
$ws.Range("E18").Value = "  +4.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.75%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.58%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.74%  "

$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.92%  "

$ws.Range("D25").Value = "3.185.69"
$ws.Range("E25").Value = "  +4.28%  "

$ws.Range("E26").Value = "  +6.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.26%  "

$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "511.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.35%  "

$ws.Range("E34").Value = "  +8.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.384"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.06%  "

$ws.Range("E40").Value = "  +6.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "188.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.84%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.91%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.49"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.96%  "

$ws.Range("E46").Value = "  +7.58%  "

$ws.Range("E47").Value = "  +2.52%  "

$ws.Range("E48").Value = "  +6.27%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.720"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.37%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.608"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.86%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.67%  "
